$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.682.76'
$ws.Range("D3").Value = '1.635.52'
$ws.Range("E3").Value = '  +1.99%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.98'
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("E6").Value = '  +2.05%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  +1.61%  '
$ws.Range("E9").Value = '  +1.70%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.00'
$ws.Range("E10").Value = '  +3.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0837'
$ws.Range("E11").Value = '  +2.67%  '
$ws.Range("D12").Value = '1.863.22'
$ws.Range("E12").Value = '  +1.91%  '
$ws.Range("D13").Value = '1.642.63'
$ws.Range("E13").Value = '  +2.36%  '
$ws.Range("E14").Value = '  +1.37%  '
$ws.Range("E15").Value = '  +2.57%  '
$ws.Range("D16").Value = '26.676.56'
$ws.Range("E16").Value = '  +1.86%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.04'
$ws.Range("E17").Value = '  +1.83%  '
$ws.Range("E18").Value = '  +1.77%  '
$ws.Range("E19").Value = '  +0.01%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '208.69'
$ws.Range("E20").Value = '  +4.00%  '
$ws.Range("E21").Value = '  +0.83%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.39'
$ws.Range("E22").Value = '  +1.04%  '
$ws.Range("E23").Value = '  +2.91%  '
$ws.Range("E24").Value = '  +1.75%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.42'
$ws.Range("E25").Value = '  +1.74%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("E27").Value = '  -0.67%  '
$ws.Range("E28").Value = '  +2.77%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.38'
$ws.Range("E29").Value = '  +1.22%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0521'
$ws.Range("E30").Value = '  +6.30%  '
$ws.Range("E31").Value = '  -0.24%  '
$ws.Range("E32").Value = '  +1.12%  '
$ws.Range("E33").Value = '  +0.82%  '
$ws.Range("E34").Value = '  +1.50%  '
$ws.Range("E35").Value = '  +0.71%  '
$ws.Range("D36").Value = '1.168.71'
$ws.Range("E36").Value = '  +0.64%  '
$ws.Range("E38").Value = '  +2.81%  '
$ws.Range("E39").Value = '  -0.02%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.32'
$ws.Range("E40").Value = '  +0.26%  '
$ws.Range("B41").Value = 'ImmutableX'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.503'
$ws.Range("E41").Value = '  +1.35%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.794'
$ws.Range("E42").Value = '  +1.50%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.37'
$ws.Range("E43").Value = '  +0.86%  '
$ws.Range("D44").Value = '1.775.50'
$ws.Range("E44").Value = '  +2.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.36'
$ws.Range("E45").Value = '  +0.83%  '
$ws.Range("E46").Value = '  +1.60%  '
$ws.Range("E47").Value = '  -1.87%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.63'
$ws.Range("E48").Value = '  +1.07%  '
$ws.Range("E49").Value = '  +1.63%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.410'
$ws.Range("E50").Value = '  +0.69%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.53'
$ws.Range("E51").Value = '  +4.13%  '
